# Scheduled-runner style price/profit refresh across leve-crafting sheets.
# Updates the computed price & profit columns (H:N) for a handful of rows
# on several job sheets; values come from an external market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 255.40741
$ws.Range("I18").Value = 179.84
$ws.Range("K18").Value = 179.84
$ws.Range("M18").Value = 104.16
$ws.Range("H43").Value = 69264.31
$ws.Range("I43").Value = 33666.332
$ws.Range("J43").Value = 77479.234
$ws.Range("K43").Value = 33666.332
$ws.Range("L43").Value = 77479.234
$ws.Range("M43").Value = -33597.332
$ws.Range("N43").Value = -77617.234
$ws.Range("H107").Value = 391.8
$ws.Range("I107").Value = 278.63635
$ws.Range("J107").Value = 703
$ws.Range("K107").Value = 278.63635
$ws.Range("L107").Value = 703
$ws.Range("M107").Value = 1641.36365
$ws.Range("N107").Value = -4543
$ws.Range("H132").Value = 217367.06
$ws.Range("I132").Value = 4911.4873
$ws.Range("K132").Value = 14734.4619
$ws.Range("M132").Value = -12204.4619
$ws.Range("H137").Value = 6283.857
$ws.Range("I137").Value = 1025.5
$ws.Range("J137").Value = 7521.1177
$ws.Range("K137").Value = 3076.5
$ws.Range("L137").Value = 22563.3531
$ws.Range("M137").Value = -526.5
$ws.Range("N137").Value = -27663.3531
$ws.Range("H141").Value = 1477
$ws.Range("I141").Value = 1505.6666
$ws.Range("J141").Value = 1305
$ws.Range("K141").Value = 4516.9998
$ws.Range("L141").Value = 3915
$ws.Range("M141").Value = 663.0002000000004
$ws.Range("N141").Value = -14275
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 11714.286
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 12000
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 12000
$ws.Range("M55").Value = -9685
$ws.Range("N55").Value = -12630
$ws.Range("H61").Value = 1327.625
$ws.Range("I61").Value = 1222.5
$ws.Range("J61").Value = 1853.25
$ws.Range("K61").Value = 1222.5
$ws.Range("L61").Value = 1853.25
$ws.Range("M61").Value = -1010.5
$ws.Range("N61").Value = -2277.25
$ws.Range("H74").Value = 1140.1224
$ws.Range("I74").Value = 1122.6562
$ws.Range("J74").Value = 1173
$ws.Range("K74").Value = 1122.6562
$ws.Range("L74").Value = 1173
$ws.Range("M74").Value = -248.6561999999999
$ws.Range("N74").Value = -2921
$ws.Range("H77").Value = 1140.1224
$ws.Range("I77").Value = 1122.6562
$ws.Range("J77").Value = 1173
$ws.Range("K77").Value = 5613.280999999999
$ws.Range("L77").Value = 5865
$ws.Range("M77").Value = -1245.280999999999
$ws.Range("N77").Value = -14601
$ws.Range("H80").Value = 17888.889
$ws.Range("J80").Value = 17888.889
$ws.Range("L80").Value = 17888.889
$ws.Range("N80").Value = -19884.889
$ws.Range("H83").Value = 17888.889
$ws.Range("J83").Value = 17888.889
$ws.Range("L83").Value = 53666.667
$ws.Range("N83").Value = -63650.667
$ws.Range("N102").ClearContents()
$ws.Range("H102").Value = 2273.9
$ws.Range("I102").Value = 2273.9
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2273.9
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -651.9000000000001
$ws.Range("H110").Value = 840.1667
$ws.Range("I110").Value = 712.2
$ws.Range("J110").Value = 1480
$ws.Range("K110").Value = 712.2
$ws.Range("L110").Value = 1480
$ws.Range("M110").Value = 1332.8
$ws.Range("N110").Value = -5570
$ws.Range("H122").Value = 1270.8334
$ws.Range("I122").Value = 1166.6666
$ws.Range("J122").Value = 1583.3334
$ws.Range("K122").Value = 3499.9998
$ws.Range("L122").Value = 4750.0002
$ws.Range("M122").Value = -1049.9998
$ws.Range("N122").Value = -9650.0002
$ws.Range("H132").Value = 171994.8
$ws.Range("I132").Value = 6322.125
$ws.Range("J132").Value = 503340.16
$ws.Range("K132").Value = 18966.375
$ws.Range("L132").Value = 1510020.48
$ws.Range("M132").Value = -16436.375
$ws.Range("N132").Value = -1515080.48
$ws.Range("H136").Value = 1327.625
$ws.Range("I136").Value = 1222.5
$ws.Range("J136").Value = 1853.25
$ws.Range("K136").Value = 3667.5
$ws.Range("L136").Value = 5559.75
$ws.Range("M136").Value = -1117.5
$ws.Range("N136").Value = -10659.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 4119.8
$ws.Range("I22").Value = 4119.8
$ws.Range("K22").Value = 4119.8
$ws.Range("M22").Value = -3946.8
$ws.Range("H134").Value = 32292454
$ws.Range("I134").Value = 2110.2354
$ws.Range("J134").Value = 71502150
$ws.Range("K134").Value = 6330.706200000001
$ws.Range("L134").Value = 214506450
$ws.Range("M134").Value = -3795.706200000001
$ws.Range("N134").Value = -214511520
$ws.Range("H141").Value = 47763.332
$ws.Range("J141").Value = 47763.332
$ws.Range("L141").Value = 47763.332
$ws.Range("N141").Value = -58123.332
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 10238
$ws.Range("J50").Value = 10238
$ws.Range("L50").Value = 10238
$ws.Range("N50").Value = -11488
$ws.Range("H51").Value = 8422
$ws.Range("J51").Value = 9371.143
$ws.Range("L51").Value = 9371.143
$ws.Range("N51").Value = -10843.143
$ws.Range("H58").Value = 1156.5938
$ws.Range("I58").Value = 994.8570999999999
$ws.Range("J58").Value = 1465.3636
$ws.Range("K58").Value = 994.8570999999999
$ws.Range("L58").Value = 1465.3636
$ws.Range("M58").Value = -791.8570999999999
$ws.Range("N58").Value = -1871.3636
$ws.Range("H61").Value = 8422
$ws.Range("J61").Value = 9371.143
$ws.Range("L61").Value = 9371.143
$ws.Range("N61").Value = -10067.143
$ws.Range("H93").Value = 6214.25
$ws.Range("I93").Value = 4295.2
$ws.Range("J93").Value = 35000
$ws.Range("K93").Value = 4295.2
$ws.Range("L93").Value = 35000
$ws.Range("M93").Value = -2423.2
$ws.Range("N93").Value = -38744
$ws.Range("H136").Value = 1156.5938
$ws.Range("I136").Value = 994.8570999999999
$ws.Range("J136").Value = 1465.3636
$ws.Range("K136").Value = 2984.5713
$ws.Range("L136").Value = 4396.0908
$ws.Range("M136").Value = -434.5712999999996
$ws.Range("N136").Value = -9496.0908
$ws.Range("H141").Value = 41626.848
$ws.Range("J141").Value = 44651.38
$ws.Range("L141").Value = 44651.38
$ws.Range("N141").Value = -55011.38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 421255.6
$ws.Range("I4").Value = 841806.4399999999
$ws.Range("J4").Value = 704.75
$ws.Range("K4").Value = 2525419.32
$ws.Range("L4").Value = 2114.25
$ws.Range("M4").Value = -2525307.32
$ws.Range("N4").Value = -2338.25
$ws.Range("H131").Value = 765.12
$ws.Range("I131").Value = 422.22223
$ws.Range("J131").Value = 840.39026
$ws.Range("K131").Value = 1266.66669
$ws.Range("L131").Value = 2521.17078
$ws.Range("M131").Value = 3773.33331
$ws.Range("N131").Value = -12601.17078
$ws.Range("H132").Value = 2131.3076
$ws.Range("J132").Value = 2131.3076
$ws.Range("L132").Value = 19181.7684
$ws.Range("N132").Value = -24241.7684
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 25390.285
$ws.Range("I102").Value = 5526.222
$ws.Range("K102").Value = 5526.222
$ws.Range("M102").Value = -3904.222
$ws.Range("H107").Value = 6376.9375
$ws.Range("I107").Value = 142.71428
$ws.Range("J107").Value = 11225.777
$ws.Range("K107").Value = 142.71428
$ws.Range("L107").Value = 11225.777
$ws.Range("M107").Value = 1777.28572
$ws.Range("N107").Value = -15065.777
$ws.Range("H132").Value = 79588
$ws.Range("I132").Value = 2663.1
$ws.Range("J132").Value = 336004.34
$ws.Range("K132").Value = 7989.299999999999
$ws.Range("L132").Value = 1008013.02
$ws.Range("M132").Value = -5459.299999999999
$ws.Range("N132").Value = -1013073.02
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 392225.12
$ws.Range("I132").Value = 130059.69
$ws.Range("J132").Value = 671868.25
$ws.Range("K132").Value = 390179.07
$ws.Range("L132").Value = 2015604.75
$ws.Range("M132").Value = -387649.07
$ws.Range("N132").Value = -2020664.75
$ws.Range("H136").Value = 372281.62
$ws.Range("I136").Value = 501674.7
$ws.Range("J136").Value = 2587.1428
$ws.Range("K136").Value = 1505024.1
$ws.Range("L136").Value = 7761.428400000001
$ws.Range("M136").Value = -1502474.1
$ws.Range("N136").Value = -12861.4284
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1333.9231
$ws.Range("I122").Value = 903.1429000000001
$ws.Range("J122").Value = 2430.4546
$ws.Range("K122").Value = 2709.4287
$ws.Range("L122").Value = 7291.3638
$ws.Range("M122").Value = -259.4287000000004
$ws.Range("N122").Value = -12191.3638
$ws.Range("H132").Value = 2601.238
$ws.Range("I132").Value = 671.1142599999999
$ws.Range("J132").Value = 5013.893
$ws.Range("K132").Value = 2013.34278
$ws.Range("L132").Value = 15041.679
$ws.Range("M132").Value = 516.6572200000001
$ws.Range("N132").Value = -20101.679
